$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 6 formatted rows below row 72 (inherits B/C formatting from the row
# above, matching style index 1), then drop the stray blank row this pushes
# down to row 79 so the sheet ends at row 78 again.
$ws.Rows("73:78").Insert(-4121)
$ws.Rows("79:79").Delete()

# Column A (Key) for the new ItemGet rows, entered first as a block.
$ws.Range("A73").Value = "ItemGet.Sword"
$ws.Range("A74").Value = "ItemGet.Boomerang"
$ws.Range("A75").Value = "ItemGet.Bow"
$ws.Range("A76").Value = "ItemGet.Fireball"
$ws.Range("A77").Value = "ItemGet.Already"

# Column B (English/Default) for sword, boomerang, bow.
$ws.Range("B73").Value = "You got the sword!"
$ws.Range("B74").Value = "You got the boomerang!"
$ws.Range("B75").Value = "You got the bow!"

# Column C (Portuguese) for sword, boomerang, bow.
$ws.Range("C73").Value = "Você obteve a espada!"
$ws.Range("C74").Value = "Você obteve o bumerangue!"
$ws.Range("C75").Value = "Você obteve o arco!"

# Fireball row: English then Portuguese.
$ws.Range("B76").Value = "You got the fireball spell!"
$ws.Range("C76").Value = "Você obteve a magia bola de fogo!"

# Already-had row: Portuguese then English.
$ws.Range("C77").Value = "Você já tinha esta arma!"
$ws.Range("B77").Value = "You already had this weapon!"

# Unknown item row, entered as a full row.
$ws.Range("A78").Value = "ItemGet.Unknown"
$ws.Range("B78").Value = "You got an unknown item!"
$ws.Range("C78").Value = "Você obteve um item desconhecido!"

$ws.Range("C66").Select()
